# Sliding window results - window 5
# Commit: Add weight setting method, implement additional activation
# functions, and enhance gradient descent optimizer.
#
# Effect on this worksheet: the model weight(s) feeding this window were
# reset to 0, so the predicted output (IPC PO, column C) collapses to 0
# for every data point. Consequently the error columns are recomputed
# from the actual/target value (IPC RO, column B):
#   C (IPC PO)   = 0
#   D (DELTA)    = C - B = -B
#   E (DELTA^2)  = D^2
# The TOTAL row re-sums the DELTA and DELTA^2 columns, and the MSE row
# re-averages DELTA^2 over the 50 data points.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 51

$deltaSum = 0
$deltaSqSum = 0

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $b = $ws.Cells.Item($r, 2).Value()

    $c = 0
    $d = $c - $b
    $e = $d * $d

    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e

    $deltaSum = $deltaSum + $d
    $deltaSqSum = $deltaSqSum + $e
}

$count = $lastRow - $firstRow + 1

# TOTAL row
$ws.Cells.Item(52, 3).Value = $deltaSum
$ws.Cells.Item(52, 5).Value = $deltaSqSum

# MSE row
$ws.Cells.Item(53, 5).Value = $deltaSqSum / $count
